$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 313; existing rows 313..422 shift down to 314..423
$ws.Rows("313:313").Insert()

# Populate the newly inserted row 313 with the new data record
$ws.Range("A313").Value = 4
$ws.Range("B313").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C313").Value = "Los Lagos"
$ws.Range("D313").Value = 44985
$ws.Range("D313").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E313").Value = 10
$ws.Range("F313").Value = 100112003
$ws.Range("G313").Value = "Ajo"
$ws.Range("H313").Value = "Chino"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 250
$ws.Range("K313").Value = 20000
$ws.Range("L313").Value = 21000
$ws.Range("M313").Value = 20400
$ws.Range("N313").Value = "$/caja 10 kilos"
$ws.Range("O313").Value = "China"
$ws.Range("P313").Value = 2040
$ws.Range("Q313").Value = 10
$ws.Range("R313").Value = "Hortaliza"
